# Applies the update described by the commit "Atualizado por script em 21-12-2023 14:45"
# - Swaps the F:V (match details / odds) content between several pairs of adjacent
#   rows (the "home"/"away" row ordering for those fixtures was corrected), while
#   leaving the Indice/pais/torneio/temporada/data_partida (A:E) columns untouched.
# - Appends two brand-new fixture rows (144 and 145) at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = "F" + $rowA + ":V" + $rowA
    $rangeB = "F" + $rowB + ":V" + $rowB
    $valuesA = $ws.Range($rangeA).Value()
    $valuesB = $ws.Range($rangeB).Value()
    $ws.Range($rangeA).Value = $valuesB
    $ws.Range($rangeB).Value = $valuesA
}

# Pairs of rows whose match/odds data (columns F through V) need to be swapped.
$pairs = @(
    @(15, 16),
    @(21, 22),
    @(31, 32),
    @(36, 37),
    @(47, 48),
    @(50, 51),
    @(71, 72),
    @(96, 97),
    @(106, 107),
    @(132, 133)
)

foreach ($pair in $pairs) {
    Swap-Rows $pair[0] $pair[1]
}

# Append two new fixture rows (144 and 145), copying the formatting (cell styles)
# from the last existing row (143) so the new cells keep the same look (bold
# bordered index column, date-time formatted match-date column, etc.).
$ws.Range("A143:V143").Copy()
$ws.Range("A144:V145").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

function Set-RowValues($rowNum, $values) {
    $arr = New-Object 'object[,]' 1, $values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $rangeAddr = "A" + $rowNum + ":V" + $rowNum
    $ws.Range($rangeAddr).Value = $arr
}

$row144 = @(
    143, "turkey", "1-lig", "2023-2024", 45281.5,
    "Boluspor", 1, "Manisa FK", 0,
    2.8, "11/12/2023 18:13", 2.92, "21/12/2023 11:59",
    3.2, "11/12/2023 18:13", 3.18, "21/12/2023 11:51",
    2.48, "11/12/2023 18:13", 2.58, "21/12/2023 11:59",
    "https://www.betexplorer.com/football/turkey/1-lig/boluspor-manisa-fk/zepTBr3g/"
)
Set-RowValues 144 $row144

$row145 = @(
    144, "turkey", "1-lig", "2023-2024", 45281.625,
    "Sanliurfaspor", 1, "Corum", 1,
    3.72, "11/12/2023 04:42", 4.03, "21/12/2023 14:55",
    3.34, "11/12/2023 04:42", 3.54, "21/12/2023 14:55",
    1.97, "11/12/2023 04:42", 1.93, "21/12/2023 14:55",
    "https://www.betexplorer.com/football/turkey/1-lig/sanliurfaspor-corum-fk/hWfOCOmm/"
)
Set-RowValues 145 $row145
